$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Hoja2")

# --- Refine existing forecast amounts (actuals came in with cents) ---
$ws1.Range("C6").Value = 1398551.39
$ws1.Range("C7").Value = 1875624.74
$ws1.Range("C10").Value = 510956.57
$ws1.Range("C12").Value = 1420750.27
$ws1.Range("C15").Value = 1405759.49
$ws1.Range("B16").Value = 2386745
$ws1.Range("C19").Value = 1071711.83
$ws1.Range("C20").Value = 1571910
$ws1.Range("B23").Value = 2370270
$ws1.Range("C23").Value = 1542160.71
$ws1.Range("C24").Value = 1430854.6
$ws1.Range("C25").Value = 1149873.39

# --- Append the next forecast day (row 26), inheriting row 25's formatting ---
$ws1.Range("A25:D25").Copy()
$ws1.Range("A26:D26").PasteSpecial(-4122)
$ws1.Range("A26").Value = 43521
$ws1.Range("B26").Value = 2030848
$ws1.Range("C26").Value = 1286885
$ws1.Range("D26").Formula = '=B26+C26/Hoja2!$A$2'

# --- Grow the Tabla1 table to cover the new row ---
$lo = $ws1.ListObjects.Item("Tabla1")
$lo.Resize($ws1.Range("A1:D26"))

# --- Update the remembered selection on Hoja2 (now sat on the rate cell) ---
$ws2.Activate() | Out-Null
$ws2.Range("A2").Select() | Out-Null
$ws1.Activate() | Out-Null
